$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.523.56'
$ws.Range('E2').Value = '  +6.02%  '
$ws.Range('D3').Value = '1.812.94'
$ws.Range('E3').Value = '  +6.01%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('E4').Value = '  +0.53%  '
$ws.Range('D5').Value = "'343.73"
$ws.Range('E5').Value = '  +3.69%  '
$ws.Range('D6').Value = "'0.9996"
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = "'0.3836"
$ws.Range('E7').Value = '  +4.03%  '
$ws.Range('D8').Value = "'50.12"
$ws.Range('E8').Value = '  +3.46%  '
$ws.Range('D9').Value = "'0.3518"
$ws.Range('E9').Value = '  +6.25%  '
$ws.Range('E10').Value = '  +4.93%  '
$ws.Range('D11').Value = "'0.07759"
$ws.Range('E11').Value = '  +4.78%  '
$ws.Range('D12').Value = "'1.002"
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('E13').Value = '  +11.88%  '
$ws.Range('D14').Value = "'6.616"
$ws.Range('E14').Value = '  +6.35%  '
$ws.Range('D15').Value = '1.814.08'
$ws.Range('E15').Value = '  +6.44%  '
$ws.Range('D16').Value = "'7.209"
$ws.Range('E16').Value = '  +4.65%  '
$ws.Range('D17').Value = "'0.00001127"
$ws.Range('E17').Value = '  +5.27%  '
$ws.Range('D18').Value = "'0.06719"
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').Value = "'86.49"
$ws.Range('E19').Value = '  +6.22%  '
$ws.Range('D20').Value = "'1.000"
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('E21').Value = '  +9.49%  '
$ws.Range('D22').Value = "'6.530"
$ws.Range('E22').Value = '  +7.73%  '
$ws.Range('D23').Value = "'13.14"
$ws.Range('E23').Value = '  +0.91%  '
$ws.Range('D24').Value = '27.526.48'
$ws.Range('E24').Value = '  +6.51%  '
$ws.Range('D25').Value = "'2.475"
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('D26').Value = "'2.678"
$ws.Range('E26').Value = '  +7.81%  '
$ws.Range('D27').Value = "'22.16"
$ws.Range('E27').Value = '  +15.46%  '
$ws.Range('D28').Value = "'1.491"
$ws.Range('E28').Value = '  +14.84%  '
$ws.Range('D29').Value = "'153.85"
$ws.Range('E29').Value = '  +2.80%  '
$ws.Range('D30').Value = '2.017.99'
$ws.Range('E30').Value = '  +6.81%  '
$ws.Range('D31').Value = "'136.71"
$ws.Range('E31').Value = '  +6.58%  '
$ws.Range('D32').Value = "'6.378"
$ws.Range('E32').Value = '  +7.11%  '
$ws.Range('D33').Value = "'4.085"
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('D34').Value = "'13.95"
$ws.Range('E34').Value = '  +7.81%  '
$ws.Range('D35').Value = "'0.08825"
$ws.Range('E35').Value = '  +3.92%  '
$ws.Range('D36').Value = "'1.723"
$ws.Range('E36').Value = '  +0.33%  '
$ws.Range('D37').Value = "'5.643"
$ws.Range('E37').Value = '  +5.63%  '
$ws.Range('D38').Value = "'0.7111"
$ws.Range('E38').Value = '  +15.58%  '
$ws.Range('D39').Value = "'0.06550"
$ws.Range('E39').Value = '  +5.47%  '
$ws.Range('D40').Value = "'0.2266"
$ws.Range('E40').Value = '  +6.83%  '
$ws.Range('D41').Value = "'0.02409"
$ws.Range('E41').Value = '  +6.13%  '
$ws.Range('D42').Value = "'9.023"
$ws.Range('E42').Value = '  +5.64%  '
$ws.Range('D43').Value = "'1.283"
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').Value = "'14.94"
$ws.Range('E44').Value = '  +2.84%  '
$ws.Range('D45').Value = "'0.6634"
$ws.Range('E45').Value = '  +13.42%  '
$ws.Range('D46').Value = "'0.9995"
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').Value = "'4.038"
$ws.Range('E47').Value = '  +5.13%  '
$ws.Range('D48').Value = "'2.183"
$ws.Range('E48').Value = '  +8.89%  '
$ws.Range('D49').Value = "'133.20"
$ws.Range('E49').Value = '  +4.94%  '
$ws.Range('E50').Value = '  +1.98%  '
$ws.Range('D51').Value = "'80.75"
$ws.Range('E51').Value = '  +5.41%  '
